$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill the remaining "x" marks in row 9 (columns C, E, G, H, I, J, K, L)
$ws.Range("C9").Value = "x"
$ws.Range("E9").Value = "x"
$ws.Range("G9").Value = "x"
$ws.Range("H9").Value = "x"
$ws.Range("I9").Value = "x"
$ws.Range("J9").Value = "x"
$ws.Range("K9").Value = "x"
$ws.Range("L9").Value = "x"

# Update the active selection to K14
$ws.Range("K14").Select()
